$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: @jobs_internships_group
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = "2026-02-12T18:31:10.067850+00:00"
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("L27").Value = "[63]"
$ws.Range("M27").Value = "[44]"

# Row 28: @oadiscussionpoint
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = "2026-02-12T18:31:27.453512+00:00"
$ws.Range("H28").Value = 3
$ws.Range("L28").Value = "[19733, 19735, 19728]"
